# Applies the LOM3004.xlsx edit:
#  - Row 13 (a spacer row holding only the "5840793 - Sérgio Schneider"
#    value in B/C, with no label in A) is removed, shifting every row
#    below it up by one.
#  - After the shift, several B/C cells end up holding the wrong
#    (now-orphaned) long-form text; those are corrected to their final
#    values below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old spacer row 13 (B13/C13 = "5840793 - Sérgio Schneider",
# no A13 label) -- everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# --- Fix up content left behind by the shift -----------------------------

# Row 10 (Objetivos:) - was the long objectives paragraph, now the
# professor's line.
$ws.Range("B10").Value = "5840793 - Sérgio Schneider"
$ws.Range("C10").Value = "5840793 - Sérgio Schneider"

# Row 13 (Programa resumido:) - was the long summary program text, now
# "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) - was the long program text, now the activation date.
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"

# Row 18 (Método:) - was "Para compor a Nota...", now the professor line.
$ws.Range("B18").Value = "5840793 - Sérgio Schneider"
$ws.Range("C18").Value = "5840793 - Sérgio Schneider"

# Row 19 (Critério:) - now holds "Para compor a Nota no Semestre..."
$ws.Range("B19").Value = "Para compor a Nota no Semestre (NS) serão feitas duas avaliações (P1 e P2)."
$ws.Range("C19").Value = "Para compor a Nota no Semestre (NS) serão feitas duas avaliações (P1 e P2)."

# Row 20 (Norma de recuperação:) - now holds the NS formula/criteria text.
$ws.Range("B20").Value = "NS = (P1 + P2)/2Serão considerados aprovados os alunos que obtiverem: NS maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Range("C20").Value = "NS = (P1 + P2)/2Serão considerados aprovados os alunos que obtiverem: NS maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."

# Row 21 (Bibliografia:) - was the long bibliography list, now the
# recovery-exam text.
$ws.Range("B21").Value = "A prova de Recuperação (R) irá compor a nota final (NF) da seguinte forma:NF = (R + NS)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Range("C21").Value = "A prova de Recuperação (R) irá compor a nota final (NF) da seguinte forma:NF = (R + NS)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
